# Update NATMI TPM-derived LR-pair metrics (Tnfsf13-Tnfrsf1a, YoungD7) with
# refreshed values following the "update scripts wuth new tpm" re-run.
# Ligand-expressing-cell counts/rates (and everything downstream of them)
# changed for the "Tnfsf13" ligand rows, and receptor-side aggregate
# columns (M:T) changed for every row whose Target cluster's underlying
# Tnfrsf1a values were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.342908333333333
$ws.Cells.Item(2, 8).Value = 4.028725
$ws.Cells.Item(2, 9).Value = 0.2879023314891748
$ws.Cells.Item(2, 10).Value = 0.2879023314891748
$ws.Cells.Item(2, 13).Value = 47.32925566666668
$ws.Cells.Item(2, 14).Value = 141.987767
$ws.Cells.Item(2, 15).Value = 0.3408416299313156
$ws.Cells.Item(2, 16).Value = 0.3408416299313156
$ws.Cells.Item(2, 17).Value = 63.55885184523056
$ws.Cells.Item(2, 18).Value = 572.029666607075
$ws.Cells.Item(2, 19).Value = 0.09812909992579624
$ws.Cells.Item(2, 20).Value = 0.09812909992579626

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.342908333333333
$ws.Cells.Item(3, 8).Value = 4.028725
$ws.Cells.Item(3, 9).Value = 0.2879023314891748
$ws.Cells.Item(3, 10).Value = 0.2879023314891748
$ws.Cells.Item(3, 13).Value = 43.717953
$ws.Cells.Item(3, 14).Value = 131.153859
$ws.Cells.Item(3, 15).Value = 0.3148348341399153
$ws.Cells.Item(3, 16).Value = 0.3148348341399154
$ws.Cells.Item(3, 17).Value = 58.709203399975
$ws.Cells.Item(3, 18).Value = 528.382830599775
$ws.Cells.Item(3, 19).Value = 0.09064168278288924
$ws.Cells.Item(3, 20).Value = 0.09064168278288927

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.342908333333333
$ws.Cells.Item(4, 8).Value = 4.028725
$ws.Cells.Item(4, 9).Value = 0.2879023314891748
$ws.Cells.Item(4, 10).Value = 0.2879023314891748
$ws.Cells.Item(4, 13).Value = 21.09134933333333
$ws.Cells.Item(4, 14).Value = 63.274048
$ws.Cells.Item(4, 15).Value = 0.1518893501062827
$ws.Cells.Item(4, 16).Value = 0.1518893501062827
$ws.Cells.Item(4, 17).Value = 28.32374878097777
$ws.Cells.Item(4, 18).Value = 254.9137390288
$ws.Cells.Item(4, 19).Value = 0.04372929802397433
$ws.Cells.Item(4, 20).Value = 0.04372929802397434

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.342908333333333
$ws.Cells.Item(5, 8).Value = 4.028725
$ws.Cells.Item(5, 9).Value = 0.2879023314891748
$ws.Cells.Item(5, 10).Value = 0.2879023314891748
$ws.Cells.Item(5, 13).Value = 26.72140366666666
$ws.Cells.Item(5, 14).Value = 80.164211
$ws.Cells.Item(5, 15).Value = 0.1924341858224863
$ws.Cells.Item(5, 16).Value = 0.1924341858224864
$ws.Cells.Item(5, 17).Value = 35.88439566233055
$ws.Cells.Item(5, 18).Value = 322.959560960975
$ws.Cells.Item(5, 19).Value = 0.05540225075651491
$ws.Cells.Item(5, 20).Value = 0.05540225075651493

$ws.Cells.Item(6, 7).Value = 0.4963216666666666
$ws.Cells.Item(6, 9).Value = 0.1064050028249084
$ws.Cells.Item(6, 10).Value = 0.1064050028249084
$ws.Cells.Item(6, 13).Value = 47.32925566666668
$ws.Cells.Item(6, 14).Value = 141.987767
$ws.Cells.Item(6, 15).Value = 0.3408416299313156
$ws.Cells.Item(6, 16).Value = 0.3408416299313156
$ws.Cells.Item(6, 17).Value = 23.49053505457278
$ws.Cells.Item(6, 18).Value = 211.414815491155
$ws.Cells.Item(6, 19).Value = 0.03626725459568802
$ws.Cells.Item(6, 20).Value = 0.03626725459568803

$ws.Cells.Item(7, 7).Value = 0.4963216666666666
$ws.Cells.Item(7, 9).Value = 0.1064050028249084
$ws.Cells.Item(7, 10).Value = 0.1064050028249084
$ws.Cells.Item(7, 13).Value = 43.717953
$ws.Cells.Item(7, 14).Value = 131.153859
$ws.Cells.Item(7, 15).Value = 0.3148348341399153
$ws.Cells.Item(7, 16).Value = 0.3148348341399154
$ws.Cells.Item(7, 19).Value = 0.03350000141603725
$ws.Cells.Item(7, 20).Value = 0.03350000141603726

$ws.Cells.Item(8, 7).Value = 0.4963216666666666
$ws.Cells.Item(8, 9).Value = 0.1064050028249084
$ws.Cells.Item(8, 10).Value = 0.1064050028249084
$ws.Cells.Item(8, 13).Value = 21.09134933333333
$ws.Cells.Item(8, 14).Value = 63.274048
$ws.Cells.Item(8, 15).Value = 0.1518893501062827
$ws.Cells.Item(8, 16).Value = 0.1518893501062827
$ws.Cells.Item(8, 17).Value = 10.46809365336889
$ws.Cells.Item(8, 18).Value = 94.21284288032
$ws.Cells.Item(8, 19).Value = 0.01616178672713251
$ws.Cells.Item(8, 20).Value = 0.01616178672713251

$ws.Cells.Item(9, 7).Value = 0.4963216666666666
$ws.Cells.Item(9, 9).Value = 0.1064050028249084
$ws.Cells.Item(9, 10).Value = 0.1064050028249084
$ws.Cells.Item(9, 13).Value = 26.72140366666666
$ws.Cells.Item(9, 14).Value = 80.164211
$ws.Cells.Item(9, 15).Value = 0.1924341858224863
$ws.Cells.Item(9, 16).Value = 0.1924341858224864
$ws.Cells.Item(9, 17).Value = 13.26241160351277
$ws.Cells.Item(9, 18).Value = 119.361704431615
$ws.Cells.Item(9, 19).Value = 0.02047596008605061
$ws.Cells.Item(9, 20).Value = 0.02047596008605061

$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1501973333333333
$ws.Cells.Item(10, 8).Value = 0.450592
$ws.Cells.Item(10, 9).Value = 0.03220038283833477
$ws.Cells.Item(10, 10).Value = 0.03220038283833477
$ws.Cells.Item(10, 13).Value = 47.32925566666668
$ws.Cells.Item(10, 14).Value = 141.987767
$ws.Cells.Item(10, 15).Value = 0.3408416299313156
$ws.Cells.Item(10, 16).Value = 0.3408416299313156
$ws.Cells.Item(10, 17).Value = 7.10872798978489
$ws.Cells.Item(10, 18).Value = 63.97855190806401
$ws.Cells.Item(10, 19).Value = 0.01097523097103038
$ws.Cells.Item(10, 20).Value = 0.01097523097103039

$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.1501973333333333
$ws.Cells.Item(11, 8).Value = 0.450592
$ws.Cells.Item(11, 9).Value = 0.03220038283833477
$ws.Cells.Item(11, 10).Value = 0.03220038283833477
$ws.Cells.Item(11, 13).Value = 43.717953
$ws.Cells.Item(11, 14).Value = 131.153859
$ws.Cells.Item(11, 15).Value = 0.3148348341399153
$ws.Cells.Item(11, 16).Value = 0.3148348341399154
$ws.Cells.Item(11, 17).Value = 6.566319959392
$ws.Cells.Item(11, 18).Value = 59.096879634528
$ws.Cells.Item(11, 19).Value = 0.0101378021901489
$ws.Cells.Item(11, 20).Value = 0.0101378021901489

$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.1501973333333333
$ws.Cells.Item(12, 8).Value = 0.450592
$ws.Cells.Item(12, 9).Value = 0.03220038283833477
$ws.Cells.Item(12, 10).Value = 0.03220038283833477
$ws.Cells.Item(12, 13).Value = 21.09134933333333
$ws.Cells.Item(12, 14).Value = 63.274048
$ws.Cells.Item(12, 15).Value = 0.1518893501062827
$ws.Cells.Item(12, 16).Value = 0.1518893501062827
$ws.Cells.Item(12, 17).Value = 3.167864426268444
$ws.Cells.Item(12, 18).Value = 28.510779836416
$ws.Cells.Item(12, 19).Value = 0.004890895222488167
$ws.Cells.Item(12, 20).Value = 0.004890895222488168

$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.1501973333333333
$ws.Cells.Item(13, 8).Value = 0.450592
$ws.Cells.Item(13, 9).Value = 0.03220038283833477
$ws.Cells.Item(13, 10).Value = 0.03220038283833477
$ws.Cells.Item(13, 13).Value = 26.72140366666666
$ws.Cells.Item(13, 14).Value = 80.164211
$ws.Cells.Item(13, 15).Value = 0.1924341858224863
$ws.Cells.Item(13, 16).Value = 0.1924341858224864
$ws.Cells.Item(13, 17).Value = 4.013483573656888
$ws.Cells.Item(13, 18).Value = 36.121352162912
$ws.Cells.Item(13, 19).Value = 0.006196454454667312
$ws.Cells.Item(13, 20).Value = 0.006196454454667314

$ws.Cells.Item(14, 7).Value = 2.675030666666667
$ws.Cells.Item(14, 8).Value = 8.025092
$ws.Cells.Item(14, 9).Value = 0.573492282847582
$ws.Cells.Item(14, 10).Value = 0.573492282847582
$ws.Cells.Item(14, 13).Value = 47.32925566666668
$ws.Cells.Item(14, 14).Value = 141.987767
$ws.Cells.Item(14, 15).Value = 0.3408416299313156
$ws.Cells.Item(14, 16).Value = 0.3408416299313156
$ws.Cells.Item(14, 17).Value = 126.6072103388405
$ws.Cells.Item(14, 18).Value = 1139.464893049564
$ws.Cells.Item(14, 19).Value = 0.1954700444388009
$ws.Cells.Item(14, 20).Value = 0.1954700444388009

$ws.Cells.Item(15, 7).Value = 2.675030666666667
$ws.Cells.Item(15, 8).Value = 8.025092
$ws.Cells.Item(15, 9).Value = 0.573492282847582
$ws.Cells.Item(15, 10).Value = 0.573492282847582
$ws.Cells.Item(15, 13).Value = 43.717953
$ws.Cells.Item(15, 14).Value = 131.153859
$ws.Cells.Item(15, 15).Value = 0.3148348341399153
$ws.Cells.Item(15, 16).Value = 0.3148348341399154
$ws.Cells.Item(15, 17).Value = 116.946864958892
$ws.Cells.Item(15, 18).Value = 1052.521784630028
$ws.Cells.Item(15, 19).Value = 0.1805553477508398
$ws.Cells.Item(15, 20).Value = 0.1805553477508399

$ws.Cells.Item(16, 7).Value = 2.675030666666667
$ws.Cells.Item(16, 8).Value = 8.025092
$ws.Cells.Item(16, 9).Value = 0.573492282847582
$ws.Cells.Item(16, 10).Value = 0.573492282847582
$ws.Cells.Item(16, 13).Value = 21.09134933333333
$ws.Cells.Item(16, 14).Value = 63.274048
$ws.Cells.Item(16, 15).Value = 0.1518893501062827
$ws.Cells.Item(16, 16).Value = 0.1518893501062827
$ws.Cells.Item(16, 17).Value = 56.42000626804622
$ws.Cells.Item(16, 18).Value = 507.7800564124161
$ws.Cells.Item(16, 19).Value = 0.0871073701326877
$ws.Cells.Item(16, 20).Value = 0.08710737013268771

$ws.Cells.Item(17, 7).Value = 2.675030666666667
$ws.Cells.Item(17, 8).Value = 8.025092
$ws.Cells.Item(17, 9).Value = 0.573492282847582
$ws.Cells.Item(17, 10).Value = 0.573492282847582
$ws.Cells.Item(17, 13).Value = 26.72140366666666
$ws.Cells.Item(17, 14).Value = 80.164211
$ws.Cells.Item(17, 15).Value = 0.1924341858224863
$ws.Cells.Item(17, 16).Value = 0.1924341858224864
$ws.Cells.Item(17, 17).Value = 71.48057426471244
$ws.Cells.Item(17, 18).Value = 643.325168382412
$ws.Cells.Item(17, 19).Value = 0.1103595205252535
$ws.Cells.Item(17, 20).Value = 0.1103595205252535

